$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q11)
$ws.Range("B7").Value = 0.188856333554674
$ws.Range("C7").Value = 0.6210211877426164
$ws.Range("D7").Value = 0.8227230598494827
$ws.Range("E7").Value = 0.9070408259000708
$ws.Range("F7").Value = 0.8997464472128907
$ws.Range("G7").Value = 36

# Row 8 (Q12)
$ws.Range("B8").Value = 0.2018765014077513
$ws.Range("C8").Value = 0.5836206403444438
$ws.Range("D8").Value = 0.7462698502533404
$ws.Range("E8").Value = 0.8638691163905214
$ws.Range("F8").Value = 0.8522125269166538
$ws.Range("G8").Value = 35

# Row 9 (Q13)
$ws.Range("B9").Value = -0.09022918269035383
$ws.Range("C9").Value = 0.3486329789925496
$ws.Range("D9").Value = 0.2017908439045041
$ws.Range("E9").Value = 0.4492113577198423
$ws.Range("F9").Value = 0.451488227386921
$ws.Range("G9").Value = 20

# Row 10 (Q14)
$ws.Range("B10").Value = -0.009659961070461246
$ws.Range("C10").Value = 0.414243314843904
$ws.Range("D10").Value = 0.2785379243470545
$ws.Range("E10").Value = 0.5277669223691975
$ws.Range("F10").Value = 0.5492252060470605
$ws.Range("G10").Value = 13

# Row 11 (Q15)
$ws.Range("B11").Value = 0.02760926664935082
$ws.Range("C11").Value = 0.3666891604559107
$ws.Range("D11").Value = 0.2015955073104189
$ws.Range("E11").Value = 0.4489938833775121
$ws.Range("F11").Value = 0.5010404620705597
$ws.Range("G11").Value = 5
